# Updating IMS assay metadata
#
# 1. dataset_type: rename "Multiplex Ion Beam Imaging" -> "MIBI" and move it
#    up in the list (new row 4, same HRAVS_0000172 term URI).
# 2. acquisition_instrument_vendor: rename "Bruker Daltonics" -> "Bruker"
#    with a new RRID (SCR_017365 instead of SCR_023608).
# 3. ms_ionization_technique: add two new ionization techniques, "HESI" and
#    "ESI", to the controlled vocabulary list.
# 4. preparation_instrument_vendor: rename "Custom" -> "In-House" (same
#    NCIt URI) and move it to the top of the list.
# 5. SIMS data sheet: widen the ms_ionization_technique dropdown validation
#    range to include the two new rows.
# 6. .metadata sheet: bump pav:createdOn to the new save timestamp.

$wb = $excel.ActiveWorkbook

# --- 1. dataset_type: "Multiplex Ion Beam Imaging" -> "MIBI" -------------
$ws = $wb.Worksheets.Item("dataset_type")
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).Value = "MIBI"
$ws.Cells.Item(4, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000172"

# --- 2. acquisition_instrument_vendor: "Bruker Daltonics" -> "Bruker" ----
$ws = $wb.Worksheets.Item("acquisition_instrument_vendor")
$ws.Cells.Item(5, 1).Value = "Bruker"
$ws.Cells.Item(5, 2).Value = "https://identifiers.org/RRID:SCR_017365"

# --- 3. ms_ionization_technique: add "HESI" and "ESI" --------------------
$ws = $wb.Worksheets.Item("ms_ionization_technique")
$ws.Rows.Item(3).Insert()
$ws.Cells.Item(3, 1).Value = "HESI"
$ws.Cells.Item(3, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000289"

$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).Value = "ESI"
$ws.Cells.Item(9, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C19363"

# --- 4. preparation_instrument_vendor: "Custom" -> "In-House", moved up --
$ws = $wb.Worksheets.Item("preparation_instrument_vendor")
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(1).Insert()
$ws.Cells.Item(1, 1).Value = "In-House"
$ws.Cells.Item(1, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C126386"

# --- 5. Widen the ms_ionization_technique validation range on SIMS -------
$ws = $wb.Worksheets.Item("SIMS")
$dv = $ws.Range("O2:O1001").Validation
$dv.Formula1 = "='ms_ionization_technique'!`$A`$1:`$A`$10"

# --- 6. .metadata: bump pav:createdOn ------------------------------------
$ws = $wb.Worksheets.Item(".metadata")
$ws.Cells.Item(2, 3).Value = "2023-10-20T20:00:51-07:00"
